$wb = $excel.ActiveWorkbook

# Work on the MainPageTest sheet (2nd sheet)
$ws = $wb.Worksheets.Item("MainPageTest")

# Add a new row of data: CheckWelcomeMessage / Welcome to Shoe Store!
$ws.Range("A5").Value = "CheckWelcomeMessage"
$ws.Range("C5").Value = "Welcome to Shoe Store!"

# Update selection on MainPageTest sheet
$ws.Range("C20").Select()

# Activate the MainPageTest sheet (making it the active tab)
$ws.Activate()

# Update selection on the MarchPageTest sheet (was previously the active tab)
$wsMarch = $wb.Worksheets.Item("MarchPageTest")
$wsMarch.Range("B17:B21").Select()

# Re-activate MainPageTest so it stays the active tab / tabSelected
$ws.Activate()
